# Finished UI Method to evaluate labview data
#
# This script:
#  1. Renames the "MP[...]" header strings to "Pmpp[...]" on all sheets
#     that use them (templateResultSheet!H1, TemplateCellGroup!H1 and
#     CellGroupsTemplate!T1 all share the same text, as do
#     CellGroupsTemplate!U1/MP_STD and V1/MP_MAX).
#  2. Moves the active sheet / selection state around: the third sheet
#     (CellGroupsTemplate) becomes the selected tab with its selection
#     moved to U5 (and the view scrolled so column O is the first
#     visible column), while the first sheet (templateResultSheet) is
#     no longer the selected tab and both it and the second sheet
#     (TemplateCellGroup) have their selected cell moved to H1.

$wb = $excel.ActiveWorkbook

$wsResult     = $wb.Worksheets.Item("templateResultSheet")
$wsCellGroup  = $wb.Worksheets.Item("TemplateCellGroup")
$wsGroupsTmpl = $wb.Worksheets.Item("CellGroupsTemplate")

# --- 1. Text fixes (shared strings) -----------------------------------
# "MP[W/cm^2]" -> "Pmpp[W/cm^2]" everywhere it is used.
$wsResult.Range("H1").Value = "Pmpp[W/cm^2]"
$wsCellGroup.Range("H1").Value = "Pmpp[W/cm^2]"
$wsGroupsTmpl.Range("T1").Value = "Pmpp[W/cm^2]"

# "MP_STD[W/cm^2]" -> "Pmpp_STD[W/cm^2]"
$wsGroupsTmpl.Range("U1").Value = "Pmpp_STD[W/cm^2]"
# "MP_MAX[W/cm^2]" -> "Pmpp_MAX[W/cm^2]"
$wsGroupsTmpl.Range("V1").Value = "Pmpp_MAX[W/cm^2]"

# --- 2. Selection / active-sheet state ---------------------------------
# First sheet: keep it not-selected-tab, move cell selection to H1.
[void]$wsResult.Activate()
[void]$wsResult.Range("H1").Select()

# Second sheet: move cell selection to H1 (tab stays unselected).
[void]$wsCellGroup.Activate()
[void]$wsCellGroup.Range("H1").Select()

# Third sheet: becomes the active tab, scrolled so column O is first
# visible, with the selection moved to U5.
[void]$wsGroupsTmpl.Activate()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
[void]$wsGroupsTmpl.Range("U5").Select()
